$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.96379729463272
$ws.Range("C2").Value = 5.7370103160159
$ws.Range("D2").Value = 13.72652201983776
$ws.Range("E2").Value = 14.0405348947244
$ws.Range("G2").Value = 3.72388950926206
$ws.Range("I2").Value = 31.80051502036775
$ws.Range("J2").Value = 8.449212012968836
$ws.Range("K2").Value = 12.16185530035756
$ws.Range("L2").Value = 12.60837900061915
$ws.Range("O2").Value = 33.32801034154452
$ws.Range("B3").Value = 14.79157065967774
$ws.Range("C3").Value = 5.660145573849797
$ws.Range("D3").Value = 13.72123069779274
$ws.Range("E3").Value = 14.06127630171764
$ws.Range("G3").Value = 3.726004812796485
$ws.Range("I3").Value = 31.89480421951064
$ws.Range("J3").Value = 8.458965415706086
$ws.Range("K3").Value = 12.04347185962412
$ws.Range("L3").Value = 12.61348765441675
$ws.Range("O3").Value = 33.41467150363764
$ws.Range("B4").Value = 14.68758061921831
$ws.Range("C4").Value = 5.611710624988963
$ws.Range("D4").Value = 13.72039603644915
$ws.Range("E4").Value = 14.07582624258129
$ws.Range("G4").Value = 3.727373024527662
$ws.Range("I4").Value = 31.95762479314491
$ws.Range("J4").Value = 8.465288192502678
$ws.Range("K4").Value = 11.97228075954733
$ws.Range("L4").Value = 12.61834408682
$ws.Range("O4").Value = 33.47350573933979
$ws.Range("B5").Value = 14.64569232744011
$ws.Range("C5").Value = 5.591671605312438
$ws.Range("D5").Value = 13.72066465557598
$ws.Range("E5").Value = 14.08221216645517
$ws.Range("G5").Value = 3.727948090178393
$ws.Range("I5").Value = 31.98446302426024
$ws.Range("J5").Value = 8.467949051749581
$ws.Range("K5").Value = 11.94367508427839
$ws.Range("L5").Value = 12.62075644218647
$ws.Range("O5").Value = 33.49889397478852
$ws.Range("B6").Value = 14.63876765423858
$ws.Range("C6").Value = 5.588326246579133
$ws.Range("D6").Value = 13.7207460697088
$ws.Range("E6").Value = 14.0833001389833
$ws.Range("G6").Value = 3.728044638475681
$ws.Range("I6").Value = 31.9889942697364
$ws.Range("J6").Value = 8.468395982906825
$ws.Range("K6").Value = 11.93895044153702
$ws.Range("L6").Value = 12.62118320725416
$ws.Range("O6").Value = 33.50319496013544
$ws.Range("B7").Value = 14.68701366038339
$ws.Range("C7").Value = 5.61144157816978
$ws.Range("D7").Value = 13.72039719233638
$ws.Range("E7").Value = 14.07591051570052
$ws.Range("G7").Value = 3.727380709100483
$ws.Range("I7").Value = 31.9579817296294
$ws.Range("J7").Value = 8.465323736213154
$ws.Range("K7").Value = 11.97189329482601
$ws.Range("L7").Value = 12.61837486509722
$ws.Range("O7").Value = 33.47384241546887
$ws.Range("B8").Value = 14.90407683282614
$ws.Range("C8").Value = 5.710770756371704
$ws.Range("D8").Value = 13.72419758524432
$ws.Range("E8").Value = 14.04731015663969
$ws.Range("G8").Value = 3.724604493908675
$ws.Range("I8").Value = 31.83200311982507
$ws.Range("J8").Value = 8.452505785705302
$ws.Range("K8").Value = 12.1207435651976
$ws.Range("L8").Value = 12.6097841025326
$ws.Range("O8").Value = 33.3567230601596
$ws.Range("B9").Value = 15.34161147481184
$ws.Range("C9").Value = 5.895282238302995
$ws.Range("D9").Value = 13.7507215410252
$ws.Range("E9").Value = 14.00560479374263
$ws.Range("G9").Value = 3.719708531232914
$ws.Range("I9").Value = 31.62407206974585
$ws.Range("J9").Value = 8.430009672818089
$ws.Range("K9").Value = 12.42321657362708
$ws.Range("L9").Value = 12.60654063370048
$ws.Range("O9").Value = 33.17173512642559
$ws.Range("B10").Value = 15.66750961381013
$ws.Range("C10").Value = 6.024022846283884
$ws.Range("D10").Value = 13.78170570245071
$ws.Range("E10").Value = 13.98370409503914
$ws.Range("G10").Value = 3.716442123642684
$ws.Range("I10").Value = 31.49517878894396
$ws.Range("J10").Value = 8.415075106352479
$ws.Range("K10").Value = 12.65012403610086
$ws.Range("L10").Value = 12.61238939931406
$ws.Range("O10").Value = 33.06314218855658
$ws.Range("B11").Value = 15.81614142925975
$ws.Range("C11").Value = 6.080998688944072
$ws.Range("D11").Value = 13.79826387375404
$ws.Range("E11").Value = 13.97563229851684
$ws.Range("G11").Value = 3.715027197734524
$ws.Range("I11").Value = 31.4417320240833
$ws.Range("J11").Value = 8.408623565183355
$ws.Range("K11").Value = 12.753986725103
$ws.Range("L11").Value = 12.61682234453704
$ws.Range("O11").Value = 33.01968713900409
$ws.Range("B12").Value = 15.87243012364965
$ws.Range("C12").Value = 6.102336650086073
$ws.Range("D12").Value = 13.80488484162091
$ws.Range("E12").Value = 13.97284699422125
$ws.Range("G12").Value = 3.714501551877343
$ws.Range("I12").Value = 31.42223962841092
$ws.Range("J12").Value = 8.406229494243252
$ws.Range("K12").Value = 12.79337728896974
$ws.Range("L12").Value = 12.6187543221597
$ws.Range("O12").Value = 33.00408771250005
$ws.Range("B13").Value = 15.86030789900808
$ws.Range("C13").Value = 6.097751843645571
$ws.Range("D13").Value = 13.80344335672954
$ws.Range("E13").Value = 13.97343480298627
$ws.Range("G13").Value = 3.71461430823762
$ws.Range("I13").Value = 31.4264044410798
$ws.Range("J13").Value = 8.406742925157809
$ws.Range("K13").Value = 12.78489166631761
$ws.Range("L13").Value = 12.61832699565212
$ws.Range("O13").Value = 33.0074092359616
$ws.Range("B14").Value = 15.82077251276282
$ws.Range("C14").Value = 6.082758987260771
$ws.Range("D14").Value = 13.79880157625422
$ws.Range("E14").Value = 13.97539771606527
$ws.Range("G14").Value = 3.714983749284379
$ws.Range("I14").Value = 31.44011340286576
$ws.Range("J14").Value = 8.408425623050618
$ws.Range("K14").Value = 12.7572263909991
$ws.Range("L14").Value = 12.61697622188751
$ws.Range("O14").Value = 33.01838660135096
$ws.Range("B15").Value = 15.79655515284578
$ws.Range("C15").Value = 6.073544231167292
$ws.Range("D15").Value = 13.79600392101373
$ws.Range("E15").Value = 13.97663537065247
$ws.Range("G15").Value = 3.715211363685373
$ws.Range("I15").Value = 31.4486078076151
$ws.Range("J15").Value = 8.409462696945289
$ws.Range("K15").Value = 12.74028747659154
$ws.Range("L15").Value = 12.61618177561151
$ws.Range("O15").Value = 33.02522207386921
$ws.Range("B16").Value = 15.65779955503881
$ws.Range("C16").Value = 6.020266583219874
$ws.Range("D16").Value = 13.78067286534763
$ws.Range("E16").Value = 13.98426960602655
$ws.Range("G16").Value = 3.716536016305845
$ws.Range("I16").Value = 31.49877611299876
$ws.Range("J16").Value = 8.415503597339169
$ws.Range("K16").Value = 12.64334651085397
$ws.Range("L16").Value = 12.6121352306163
$ws.Range("O16").Value = 33.06610180536155
$ws.Range("B17").Value = 15.57273912109548
$ws.Range("C17").Value = 5.987168929499956
$ws.Range("D17").Value = 13.77189630077313
$ws.Range("E17").Value = 13.98943686447054
$ws.Range("G17").Value = 3.717366790964811
$ws.Range("I17").Value = 31.53088193729291
$ws.Range("J17").Value = 8.419296992819826
$ws.Range("K17").Value = 12.5840176308234
$ws.Range("L17").Value = 12.61010568191874
$ws.Range("O17").Value = 33.09270363740575
$ws.Range("B18").Value = 15.52385320460521
$ws.Range("C18").Value = 5.967982995670109
$ws.Range("D18").Value = 13.76708034006231
$ws.Range("E18").Value = 13.99258695886877
$ws.Range("G18").Value = 3.71785131449659
$ws.Range("I18").Value = 31.54983657838476
$ws.Range("J18").Value = 8.421511083085161
$ws.Range("K18").Value = 12.5499555087706
$ws.Range("L18").Value = 12.60910525743329
$ws.Range("O18").Value = 33.1085637568484
$ws.Range("B19").Value = 15.50730942369251
$ws.Range("C19").Value = 5.961461653754075
$ws.Range("D19").Value = 13.76548970050409
$ws.Range("E19").Value = 13.99368412003475
$ws.Range("G19").Value = 3.718016515365396
$ws.Range("I19").Value = 31.55633812212714
$ws.Range("J19").Value = 8.422266278497473
$ws.Range("K19").Value = 12.53843434219943
$ws.Range("L19").Value = 12.60879524415017
$ws.Range("O19").Value = 33.11402976653699
$ws.Range("B20").Value = 15.58179031621337
$ws.Range("C20").Value = 5.990707723228548
$ws.Range("D20").Value = 13.77280658510802
$ws.Range("E20").Value = 13.98886838038282
$ws.Range("G20").Value = 3.717277662193099
$ws.Range("I20").Value = 31.52741368273668
$ws.Range("J20").Value = 8.418889845433664
$ws.Range("K20").Value = 12.59032707981961
$ws.Range("L20").Value = 12.61030446619176
$ws.Range("O20").Value = 33.08981391656643
$ws.Range("B21").Value = 15.83238526649518
$ws.Range("C21").Value = 6.087169269984758
$ws.Range("D21").Value = 13.80015548972302
$ws.Range("E21").Value = 13.97481380273846
$ws.Range("G21").Value = 3.714874960268876
$ws.Range("I21").Value = 31.43606647616401
$ws.Range("J21").Value = 8.407930046202445
$ws.Range("K21").Value = 12.76535097920476
$ws.Range("L21").Value = 12.61736611479141
$ws.Range("O21").Value = 33.01513904045337
$ws.Range("B22").Value = 15.9961707665
$ws.Range("C22").Value = 6.148823908733346
$ws.Range("D22").Value = 13.82007217826617
$ws.Range("E22").Value = 13.96720940521547
$ws.Range("G22").Value = 3.713363828570797
$ws.Range("I22").Value = 31.38071868342771
$ws.Range("J22").Value = 8.401052616180509
$ws.Range("K22").Value = 12.88007436186817
$ws.Range("L22").Value = 12.62345703846742
$ws.Range("O22").Value = 32.9713249494678
$ws.Range("B23").Value = 15.90877078998445
$ws.Range("C23").Value = 6.116047609646307
$ws.Range("D23").Value = 13.80925659605872
$ws.Range("E23").Value = 13.97112356052639
$ws.Range("G23").Value = 3.71416495037223
$ws.Range("I23").Value = 31.40986029846551
$ws.Range("J23").Value = 8.404697188439828
$ws.Range("K23").Value = 12.81882421801853
$ws.Range("L23").Value = 12.6200717085765
$ws.Range("O23").Value = 32.9942523612616
$ws.Range("B24").Value = 15.57769821584346
$ws.Range("C24").Value = 5.989108324965345
$ws.Range("D24").Value = 13.77239432949824
$ws.Range("E24").Value = 13.98912483328339
$ws.Range("G24").Value = 3.717317935813291
$ws.Range("I24").Value = 31.52898013384464
$ws.Range("J24").Value = 8.419073813282854
$ws.Range("K24").Value = 12.58747442955968
$ws.Range("L24").Value = 12.6102140774211
$ws.Range("O24").Value = 33.09111859520216
$ws.Range("B25").Value = 15.22226429749886
$ws.Range("C25").Value = 5.846525715572142
$ws.Range("D25").Value = 13.74151651722875
$ws.Range("E25").Value = 14.01535023656986
$ws.Range("G25").Value = 3.720974698417338
$ws.Range("I25").Value = 31.67613254170012
$ws.Range("J25").Value = 8.435814514011994
$ws.Range("K25").Value = 12.34043866649845
$ws.Range("L25").Value = 12.60596761780058
$ws.Range("O25").Value = 33.21698723409378
